# Master Acquisitions.xlsx - "Just a few changes to suppliers/notes"
#
# The 10uF capacitor (row 9 of the acquisitions table, refdes C402) is being
# re-sourced from a Murata part to a Taiyo Yuden part. Update the
# Supplier PN / Manufacturer / Manufacturer PN / Pricing columns for that row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 9 -> 10uF ceramic capacitor (C402): switch supplier part from the
# Murata GRM188R61A106KE69D to the Taiyo Yuden AMK107ABJ106MAHT.
$ws.Range("G9").Value = "963-AMK107ABJ106MAHT"
$ws.Range("I9").Value = "AMK107ABJ106MAHT"
$ws.Range("J9").Value = "0.073/0.042/--"
$ws.Range("H9").Value = "Taiyo Yuden"

# Update the saved view/selection state of the sheet.
$ws.Range("I13").Select()
